$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Add new row 9 with a Cluster_8 entry and its taxonomy
$ws.Range("A9").Value = "Cluster_8"
$ws.Range("B9").Value = "Viruses;Riboviria;Orthornavirae;Kitrinoviricota;Alsuviricetes;Tymovirales;Betaflexiviridae;Trivirinae;Chordovirus;Carrot Ch virus 1"

# Update the selected cell to match the new active cell
$ws.Range("B9").Select()
